$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Header row: uppercase WEEKDAY/WEEKEND and drop "(hours after midnight)" ---
$t.Cell(1, 2).Range.Text = "WEEKDAY earliest arrival time"
$t.Cell(1, 3).Range.Text = "WEEKDAY latest departure time"
$t.Cell(1, 4).Range.Text = "WEEKDAY highest occupant count"
$t.Cell(1, 5).Range.Text = "WEEKEND earliest arrival time"
$t.Cell(1, 6).Range.Text = "WEEKEND latest departure time"
$t.Cell(1, 7).Range.Text = "WEEKEND highest occupant count"

# --- Data rows: reformat non-empty arrival/departure time values as zero-padded HH:00 ---
# Columns 2 & 3 = weekday earliest/latest, columns 5 & 6 = weekend earliest/latest.
# Column 4 & 7 (occupant counts) and column 1 (floor) are left untouched.
for ($r = 2; $r -le $t.Rows.Count; $r++) {
    foreach ($c in 2, 3, 5, 6) {
        $cell = $t.Cell($r, $c)
        $val = $cell.Range.Text
        # strip trailing cell-end marker characters (cr + cell-mark)
        $val = $val -replace "[\x07\x0d]+$", ""
        if ($val -ne "") {
            $formatted = "{0:D2}:00" -f [int]$val
            $cell.Range.Text = $formatted
        }
    }
}
